# Add a new NPC ("Wondering Merchant") to the NPCs sheet and a matching
# row on the Npcs Commands sheet, per the commit:
#   "Added new quest items, npc and quest." (WonderingMerchant / Voidance)

$wb = $excel.ActiveWorkbook

# ----- NPCs sheet -----
$npcs = $wb.Worksheets.Item("NPCs")

$npcs.Cells.Item(9, 1).Value = "WonderingMerchant"
$npcs.Cells.Item(9, 2).Value = "Wondering Merchant"
$npcs.Cells.Item(9, 3).Value = 2
$npcs.Cells.Item(9, 4).Value = "Surface"
$npcs.Cells.Item(9, 6).Value = 1
$npcs.Cells.Item(9, 7).Value = "/m WonderingMerchant:"
$npcs.Cells.Item(9, 8).Value = 32
$npcs.Cells.Item(9, 9).Value = 256

# Widen columns A and B to fit the new, longer NPC name/text values
# (xlsx <col> width = ColumnWidth + 5/6, so subtract that padding here
# to land on exactly 21 / 22 in the saved file).
$npcs.Columns.Item(1).ColumnWidth = 21 - 5/6
$npcs.Columns.Item(2).ColumnWidth = 22 - 5/6

# ----- Npcs Commands sheet -----
$cmds = $wb.Worksheets.Item("Npcs Commands")

$cmds.Cells.Item(9, 1).Value = "Wondering Merchant"
$cmds.Cells.Item(9, 2).Value = "Voidance"
$cmds.Cells.Item(9, 3).Value = 0

# Widen column A to fit the new, longer npc_id value
$cmds.Columns.Item(1).ColumnWidth = 22 - 5/6
